$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row 7 (shifting the previously-empty row 7 and below down one
# row) copying the formatting from row 6 (the row directly above), then
# remove the now-duplicated old row (which shifted to row 8). Net effect:
# row 7 keeps the same cell styles used by the rows above it (A7 gets the
# blue/bordered style, B7/C7/H7 get the wrap-text style, G7 keeps the
# checkmark style) while every other row stays exactly where it was.
$ws.Rows.Item(7).Insert(-4121, 0)   # xlShiftDown, xlFormatFromLeftOrAbove
$ws.Rows.Item(8).Delete(-4121)      # xlShiftUp

# Fill in the new "Swap Salary" problem entry
$ws.Range("A7").Value = 627
$ws.Range("B7").Value = "Swap Salary"
$ws.Range("C7").Value = "SELECT and ORDER"
$ws.Range("D7").Value = "UPDATE"
$ws.Range("E7").Value = "Easy"
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = [char]0x2705
$ws.Range("H7").Value = "Copied 1 sol from sols and understood"

# Row grows taller (to 30) to match the wrapped-text rows above it
$ws.Rows.Item(7).RowHeight = 30

# Move the active selection to H8, matching where editing left off
$ws.Range("H8").Select()
